{"js": "// Letter edit: bump the date, split the sender's address onto two lines,\n// and drop the blank \"No Spacing\" paragraph that used to separate the\n// \"Board of Directors\" line from the following Title-styled spacer.\n\nconst body = context.document.body;\n\n// 1) \"September 19, 2025\" -> \"September 21, 2025\" (the dateline near the\n//    top of the letter; leave any other occurrences untouched).\nconst dateResults = body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"September 21, 2025\", \"Replace\");\n}\n\n// 2) Split \"3370 Eichers Pl, Santa Clara CA 95051\" (the sender address line,\n//    i.e. not the one inside the PROPERTY ADDRESS table cell) into two\n//    paragraphs: \"3370 Eichers Pl\" and \"Santa Clara, CA 95051\", each\n//    keeping the original paragraph's formatting.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"3370 Eichers Pl, Santa Clara CA 95051\") {\n    addressParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (addressParagraph) {\n  addressParagraph.insertParagraph(\"Santa Clara, CA 95051\", \"After\");\n  addressParagraph.insertText(\"3370 Eichers Pl\", \"Replace\");\n}\nawait context.sync();\n\n// 3) Remove the empty \"No Spacing\" paragraph directly after the\n//    \"Board of Directors\" signature line.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nlet boardIdx = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    boardIdx = i;\n    break;\n  }\n}\nif (boardIdx !== -1 && boardIdx + 1 < paragraphs2.items.length) {\n  const nextParagraph = paragraphs2.items[boardIdx + 1];\n  nextParagraph.load(\"text\");\n  await context.sync();\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Letter edit: bump the date, split the sender's address onto two lines,\n# and drop the blank \"No Spacing\" paragraph that used to separate the\n# \"Board of Directors\" line from the following Title-styled spacer.\n\n$d = $word.ActiveDocument\n\n# 1) \"September 19, 2025\" -> \"September 21, 2025\" (the dateline near the\n#    top of the letter; leave any other occurrences untouched).\n$dateCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $dateCount; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($ptext -eq \"September 19, 2025\") {\n    $p.Range.Text = \"September 21, 2025\"\n    break\n  }\n}\n\n# 2) Split \"3370 Eichers Pl, Santa Clara CA 95051\" (the sender address line,\n#    i.e. not the one inside the PROPERTY ADDRESS table cell) into two\n#    paragraphs: \"3370 Eichers Pl\" and \"Santa Clara, CA 95051\", each\n#    keeping the original paragraph's formatting.\n$count = $d.Paragraphs.Count\n$addressIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $ptext = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($ptext -eq \"3370 Eichers Pl, Santa Clara CA 95051\") {\n    $addressIndex = $i\n    break\n  }\n}\nif ($addressIndex -ne -1) {\n  $addressParagraph = $d.Paragraphs.Item($addressIndex)\n  $addressParagraph.Range.InsertParagraphAfter()\n  $newParagraph = $d.Paragraphs.Item($addressIndex + 1)\n  $newParagraph.Range.InsertBefore(\"Santa Clara, CA 95051\")\n  $addressParagraph.Range.Text = \"3370 Eichers Pl\"\n}\n\n# 3) Remove the empty \"No Spacing\" paragraph directly after the\n#    \"Board of Directors\" signature line.\n$count2 = $d.Paragraphs.Count\n$boardIndex = -1\nfor ($i = 1; $i -le $count2; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"*Board of Directors*\") {\n    $boardIndex = $i\n    break\n  }\n}\nif ($boardIndex -ne -1 -and ($boardIndex + 1) -le $count2) {\n  $nextParagraph = $d.Paragraphs.Item($boardIndex + 1)\n  $nextText = $nextParagraph.Range.Text.TrimEnd([char]13, [char]7)\n  if ($nextText -eq \"\") {\n    $nextParagraph.Range.Delete()\n  }\n}\n"}
